$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values as literal text (e.g. "29.30", "0.0640")
# in the source workbook, even though they look numeric. Assigning a plain
# numeric-looking string via .Value lets Excel auto-convert the cell to a
# Number (silently dropping meaningful trailing zeros, e.g. "29.30" -> 29.3).
# To keep these cells as text - matching the original file exactly - force a
# Text number format before assigning, then restore the default "Normal" style
# afterwards so no visible/structural formatting change is left behind.
# Column E ("Volume(1h)") already stays text because of its leading/trailing
# spaces, so it needs no special handling.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.055.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.681.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "220.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.533"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.266"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0639"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.920.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.682.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.606"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.09"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.993.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "66.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.64%  "
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.22%  "
$ws.Range("E31").Value = "  +4.02%  "
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.515.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.27%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.73"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "84.39"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +12.77%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.610"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.68%  "
$ws.Range("E39").Value = "  +5.42%  "
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +4.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.840"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("E45").Value = "  +2.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "51.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.812.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("E50").Value = "  +8.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "93.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.21%  "
